$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'50.930.97"
$ws.Range("E2").Value = "'  -2.76%  "
$ws.Range("D3").Value = "'2.900.05"
$ws.Range("E3").Value = "'  -2.61%  "
$ws.Range("E4").Value = "'  +0.13%  "
$ws.Range("D5").Value = "'369.37"
$ws.Range("E5").Value = "'  +3.15%  "
$ws.Range("D6").Value = "'104.13"
$ws.Range("E6").Value = "'  -4.77%  "
$ws.Range("D7").Value = "'0.541"
$ws.Range("E7").Value = "'  -5.75%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "'  -0.05%  "
$ws.Range("D9").Value = "'0.588"
$ws.Range("E9").Value = "'  -6.63%  "
$ws.Range("D10").Value = "'36.95"
$ws.Range("E10").Value = "'  -5.45%  "
$ws.Range("E11").Value = "'  +0.62%  "
$ws.Range("D12").Value = "'0.0834"
$ws.Range("E12").Value = "'  -4.39%  "
$ws.Range("D13").Value = "'18.35"
$ws.Range("E13").Value = "'  -5.93%  "
$ws.Range("D14").Value = "'3.366.90"
$ws.Range("E14").Value = "'  -1.66%  "
$ws.Range("D15").Value = "'7.35"
$ws.Range("E15").Value = "'  -5.73%  "
$ws.Range("D16").Value = "'2.904.46"
$ws.Range("E16").Value = "'  -2.50%  "
$ws.Range("D17").Value = "'0.944"
$ws.Range("E17").Value = "'  -4.08%  "
$ws.Range("D18").Value = "'50.978.20"
$ws.Range("E18").Value = "'  -2.55%  "
$ws.Range("D19").Value = "'3.29"
$ws.Range("E19").Value = "'  -5.62%  "
$ws.Range("D20").Value = "'7.25"
$ws.Range("E20").Value = "'  -5.64%  "
$ws.Range("D21").Value = "'12.96"
$ws.Range("E21").Value = "'  -7.04%  "
$ws.Range("D22").Value = "'0.0₃0942"
$ws.Range("E22").Value = "'  -4.46%  "
$ws.Range("D23").Value = "'68.25"
$ws.Range("E23").Value = "'  -3.24%  "
$ws.Range("D24").Value = "'259.30"
$ws.Range("E24").Value = "'  -4.52%  "
$ws.Range("D25").Value = "'2.70"
$ws.Range("E25").Value = "'  -3.53%  "
$ws.Range("D26").Value = "'4.34"
$ws.Range("E26").Value = "'  +3.93%  "
$ws.Range("D27").Value = "'0.171"
$ws.Range("E27").Value = "'  -4.57%  "
$ws.Range("E28").Value = "'  -0.08%  "
$ws.Range("B29").Value = "'EthereumClassic"
$ws.Range("C29").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'25.79"
$ws.Range("E29").Value = "'  -5.44%  "
$ws.Range("B30").Value = "'Filecoin"
$ws.Range("C30").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "'7.31"
$ws.Range("E30").Value = "'  -7.19%  "
$ws.Range("E31").Value = "'  -5.17%  "
$ws.Range("B32").Value = "'RenderToken"
$ws.Range("C32").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").Value = "'6.10"
$ws.Range("E32").Value = "'  -0.66%  "
$ws.Range("B33").Value = "'Cosmos"
$ws.Range("C33").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "'9.89"
$ws.Range("E33").Value = "'  -5.87%  "
$ws.Range("B34").Value = "'Toncoin"
$ws.Range("C34").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").Value = "'2.13"
$ws.Range("E34").Value = "'  -1.89%  "
$ws.Range("B35").Value = "'InjectiveProtocol"
$ws.Range("C35").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "'34.92"
$ws.Range("E35").Value = "'  -8.05%  "
$ws.Range("D36").Value = "'50.99"
$ws.Range("E36").Value = "'  -2.37%  "
$ws.Range("E37").Value = "'  +0.36%  "
$ws.Range("D38").Value = "'0.0421"
$ws.Range("E38").Value = "'  -5.13%  "
$ws.Range("D39").Value = "'3.08"
$ws.Range("E39").Value = "'  -5.04%  "
$ws.Range("D40").Value = "'2.68"
$ws.Range("E40").Value = "'  -2.02%  "
$ws.Range("D41").Value = "'16.94"
$ws.Range("E41").Value = "'  -7.14%  "
$ws.Range("D42").Value = "'1.85"
$ws.Range("E42").Value = "'  -8.63%  "
$ws.Range("E43").Value = "'  -6.22%  "
$ws.Range("D44").Value = "'22.19"
$ws.Range("E44").Value = "'  -5.53%  "
$ws.Range("D45").Value = "'117.65"
$ws.Range("E45").Value = "'  -0.85%  "
$ws.Range("E46").Value = "'  -3.71%  "
$ws.Range("D47").Value = "'2.042.45"
$ws.Range("E47").Value = "'  -4.87%  "
$ws.Range("D48").Value = "'2.32"
$ws.Range("E48").Value = "'  -5.99%  "
$ws.Range("D49").Value = "'3.17"
$ws.Range("E49").Value = "'  -8.07%  "
$ws.Range("D50").Value = "'3.205.99"
$ws.Range("E50").Value = "'  -1.33%  "
$ws.Range("E51").Value = "'  -3.23%  "